# Update "Last Updated" timestamp on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 09:45 AM"

# Refresh the "Stock List" sheet: drop the two oldest rows (MIDWESTLTD,
# CAPTRU-RE1) from the top -- the remaining rows shift up -- and append the
# two newest rows (SMARTWORKS, TRAVELFOOD) at the bottom.
$wsStock = $wb.Worksheets.Item("Stock List")
$wsStock.Range("A2:H3").EntireRow.Delete()

$wsStock.Range("A75").Value = "📋"
$wsStock.Range("B75").Value = "SMARTWORKS"
$wsStock.Range("C75").Value = "SMARTWORKS"
$wsStock.Range("D75").Value = 606.65
$wsStock.Range("E75").Value = 2.0867
$wsStock.Range("F75").Value = "N/A"
$wsStock.Range("G75").Value = "N/A"
$wsStock.Range("H75").Value = 6931.2448

$wsStock.Range("A76").Value = "📋"
$wsStock.Range("B76").Value = "TRAVELFOOD"
$wsStock.Range("C76").Value = "TRAVELFOOD"
$wsStock.Range("D76").Value = 1316.3
$wsStock.Range("E76").Value = 0.1141
$wsStock.Range("F76").Value = "N/A"
$wsStock.Range("G76").Value = "N/A"
$wsStock.Range("H76").Value = 17332.9705
